# Update workTime sheet for month 6:
#  - Insert a new column A ("登记号") holding the numeric employee id
#    (leading zeros stripped) ahead of the existing 员工工号 / 姓名 columns.
#  - The previous column A (员工工号) shifts to column B, and the previous
#    column B (姓名) shifts to column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:B to B:C, leaving a blank column A.
$ws.Columns("A").Insert()

# New header for column A.
$ws.Range("A1").Value = "登记号"

# Numeric registration numbers for rows 2-92 (derived from the old
# zero-padded 员工工号 text values, now living in column B).
$regNos = @(109, 314, 311, 511, 20022, 20025, 563, 564, 565, 566, 598, 607, 653, 654, 727, 728, 729, 781, 782, 783, 795, 849, 851, 852, 854, 855, 856, 857, 858, 950, 951, 952, 954, 956, 957, 958, 960, 962, 889, 1008, 1009, 1010, 1011, 1012, 1013, 1014, 1092, 1093, 1094, 1096, 1175, 1176, 1178, 1179, 1184, 1185, 1186, 1187, 1188, 1208, 1217, 1226, 1152, 1228, 1327, 1381, 1382, 1383, 1413, 1414, 1455, 1456, 1457, 1458, 1459, 1543, 1544, 1536, 1540, 1542, 1634, 1635, 1636, 1642, 1794, 1795, 1796, 1797, 1774, 1779, 1708)

for ($i = 0; $i -lt $regNos.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $regNos[$i]
}

# Match the author's final selection (cell I10) recorded in the sheet view.
$ws.Range("I10").Select()
